$d = $word.ActiveDocument

# Locate the paragraph containing the stand-alone "Ι.Κ.Υ." bullet item
# (a <w:p> right before the "${local_directorate}" bullet under the
# "ΚΟΙΝΟΠΟΙΗΣΗ" heading) and remove it entirely, paragraph mark included,
# so the following paragraph takes its place in the list.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -match "^\s*Ι\.Κ\.Υ\.\s*$") {
        $para.Range.Delete()
        $found = $true
        break
    }
}

Write-Output "Deleted paragraph: $found"
